# Apply updated loading_percent results for Case_3_57 (380 kV case).
# Updates numeric results in columns B, C, D, F, G, N for rows 2-25
# to reflect the new simulation output values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.00648322401483
$ws.Range("C2").Value = 10.69828057976813
$ws.Range("D2").Value = 4.981623389126059
$ws.Range("F2").Value = 29.55748212936986
$ws.Range("G2").Value = 3.632106957472635
$ws.Range("N2").Value = 17.00829988987423
$ws.Range("B3").Value = 16.35302150500121
$ws.Range("C3").Value = 10.07373862251478
$ws.Range("D3").Value = 5.015163060677339
$ws.Range("F3").Value = 29.27485066424385
$ws.Range("G3").Value = 3.635999053925908
$ws.Range("N3").Value = 17.08916981434259
$ws.Range("B4").Value = 15.94320722061095
$ws.Range("C4").Value = 9.672007426524317
$ws.Range("D4").Value = 5.036627821544934
$ws.Range("F4").Value = 29.11213250202849
$ws.Range("G4").Value = 3.638510617883818
$ws.Range("N4").Value = 17.14090347359768
$ws.Range("B5").Value = 15.77433220721966
$ws.Range("C5").Value = 9.503855165650759
$ws.Range("D5").Value = 5.045594692439061
$ws.Range("F5").Value = 29.04860873339883
$ws.Range("G5").Value = 3.639564849794331
$ws.Range("N5").Value = 17.16251028206955
$ws.Range("B6").Value = 15.74618648666673
$ws.Range("C6").Value = 9.475670185097604
$ws.Range("D6").Value = 5.047096931502469
$ws.Range("F6").Value = 29.03823063079917
$ws.Range("G6").Value = 3.639741764937086
$ws.Range("N6").Value = 17.16612984142947
$ws.Range("B7").Value = 15.94093689466966
$ws.Range("C7").Value = 9.669757439146853
$ws.Range("D7").Value = 5.036747860997988
$ws.Range("F7").Value = 29.11126444106202
$ws.Range("G7").Value = 3.638524710963195
$ws.Range("N7").Value = 17.14119274249249
$ws.Range("B8").Value = 16.78311493440999
$ws.Range("C8").Value = 10.48680788211401
$ws.Range("D8").Value = 4.993007470704613
$ws.Range("F8").Value = 29.45782638496409
$ws.Range("G8").Value = 3.633423745760905
$ws.Range("N8").Value = 17.03575356874946
$ws.Range("B9").Value = 18.35532964234293
$ws.Range("C9").Value = 11.93908014991649
$ws.Range("D9").Value = 4.914112422918913
$ws.Range("F9").Value = 30.22009429542702
$ws.Range("G9").Value = 3.624381662019557
$ws.Range("N9").Value = 16.8453878117096
$ws.Range("B10").Value = 19.44907093696457
$ws.Range("C10").Value = 12.90946038661791
$ws.Range("D10").Value = 4.860294823979975
$ws.Range("F10").Value = 30.82584676095372
$ws.Range("G10").Value = 3.61831646158958
$ws.Range("N10").Value = 16.7153849479996
$ws.Range("B11").Value = 19.9310806596658
$ws.Range("C11").Value = 13.32920626576599
$ws.Range("D11").Value = 4.836702172633514
$ws.Range("F11").Value = 31.11025059025355
$ws.Range("G11").Value = 3.615681093220704
$ws.Range("N11").Value = 16.65835407118458
$ws.Range("B12").Value = 20.11121564355559
$ws.Range("C12").Value = 13.48499173400841
$ws.Range("D12").Value = 4.827895399106762
$ws.Range("F12").Value = 31.21912190702791
$ws.Range("G12").Value = 3.614700810709469
$ws.Range("N12").Value = 16.63705884414714
$ws.Range("B13").Value = 20.07252900139363
$ws.Range("C13").Value = 13.45158182982517
$ws.Range("D13").Value = 4.829786445533124
$ws.Range("F13").Value = 31.19562390108752
$ws.Range("G13").Value = 3.614911147873526
$ws.Range("N13").Value = 16.64163179150297
$ws.Range("B14").Value = 19.94594923111199
$ws.Range("C14").Value = 13.34208644354469
$ws.Range("D14").Value = 4.835975087272577
$ws.Range("F14").Value = 31.11918454425946
$ws.Range("G14").Value = 3.6156000912367
$ws.Range("N14").Value = 16.65659607648882
$ws.Range("B15").Value = 19.8680996986742
$ws.Range("C15").Value = 13.27460427167941
$ws.Range("D15").Value = 4.839782362750373
$ws.Range("F15").Value = 31.07251314252931
$ws.Range("G15").Value = 3.61602438719473
$ws.Range("N15").Value = 16.66580128929759
$ws.Range("B16").Value = 19.41724386712814
$ws.Range("C16").Value = 12.88158855374856
$ws.Range("D16").Value = 4.861854499707798
$ws.Range("F16").Value = 30.80742994257739
$ws.Range("G16").Value = 3.618491168721719
$ws.Range("N16").Value = 16.71915429785625
$ws.Range("B17").Value = 19.13655992847566
$ws.Range("C17").Value = 12.63489651152047
$ws.Range("D17").Value = 4.875622351408937
$ws.Range("F17").Value = 30.64700620415433
$ws.Range("G17").Value = 3.620036064226988
$ws.Range("N17").Value = 16.75242313772152
$ws.Range("B18").Value = 18.97366335134805
$ws.Range("C18").Value = 12.49096657409505
$ws.Range("D18").Value = 4.883624991723049
$ws.Range("F18").Value = 30.55557319952402
$ws.Range("G18").Value = 3.620936299497772
$ws.Range("N18").Value = 16.77175703677358
$ws.Range("B19").Value = 18.91826468025555
$ws.Range("C19").Value = 12.44188569004283
$ws.Range("D19").Value = 4.886348949117463
$ws.Range("F19").Value = 30.52476251221086
$ws.Range("G19").Value = 3.621243108535336
$ws.Range("N19").Value = 16.77833732390475
$ws.Range("B20").Value = 19.16659096433131
$ws.Range("C20").Value = 12.66136871051456
$ws.Range("D20").Value = 4.874148078263476
$ws.Range("F20").Value = 30.66399750730347
$ws.Range("G20").Value = 3.619870402347602
$ws.Range("N20").Value = 16.74886107913844
$ws.Range("B21").Value = 19.98319483523797
$ws.Range("C21").Value = 13.37433404466827
$ws.Range("D21").Value = 4.834153884482625
$ws.Range("F21").Value = 31.14160557266149
$ws.Range("G21").Value = 3.615397253123882
$ws.Range("N21").Value = 16.65219254671415
$ws.Range("B22").Value = 20.50289384205021
$ws.Range("C22").Value = 13.83805691341151
$ws.Range("D22").Value = 4.808756808346468
$ws.Range("F22").Value = 31.46054480438615
$ws.Range("G22").Value = 3.612576756872961
$ws.Range("N22").Value = 16.59076834824446
$ws.Range("B23").Value = 20.22684820232422
$ws.Range("C23").Value = 13.584701333171
$ws.Range("D23").Value = 4.822244065196632
$ws.Range("F23").Value = 31.28973211329585
$ws.Range("G23").Value = 3.614072725958675
$ws.Range("N23").Value = 16.62339174589442
$ws.Range("B24").Value = 19.15301867913571
$ws.Range("C24").Value = 12.64940717744131
$ws.Range("D24").Value = 4.874814325126434
$ws.Range("F24").Value = 30.65631325007728
$ws.Range("G24").Value = 3.619945260528337
$ws.Range("N24").Value = 16.75047084029043
$ws.Range("B25").Value = 17.93994020185811
$ws.Range("C25").Value = 11.56286823933709
$ws.Range("D25").Value = 4.93472405297923
$ws.Range("F25").Value = 30.00549372921323
$ws.Range("G25").Value = 3.626725719250617
$ws.Range("N25").Value = 16.89514544084691
